$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correção das notas do fórum: zera todas as notas semanais (colunas B:Q)
# para as linhas de alunos (linhas 3 a 50), mantendo a matrícula (coluna A).
$ws.Range("B3:Q50").Value = 0
